$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 10 - the note text gains a name
$ws.Range("B10").Value = "Carmen en Steven"

# New "Klantservice" feature rows (19-26)
$ws.Range("A19").Value = "Klantservice"
$ws.Range("B19").Value = "Steven"

$ws.Range("A20").Borders.LineStyle = -4142
$ws.Range("A20").Value = "KlantDAO"
$ws.Range("B20").Value = "Steven"

$ws.Range("A21").Value = "JDBCKlantDAO"
$ws.Range("B21").Value = "Steven"

$ws.Range("A22").Value = "KlantController"
$ws.Range("B22").Value = "Steven, Wim en Elise"
$ws.Range("C22").Value = "Steven"

$ws.Range("A23").Value = "RegistratieService"
$ws.Range("B23").Value = "Steven"

$ws.Range("A24").Value = "RegistrationDto"
$ws.Range("B24").Value = "Steven"

$ws.Range("A25").Value = "Model Klant"
$ws.Range("B25").Value = "Steven"

$ws.Range("A26").Value = "Model Gebruiker"
$ws.Range("B26").Value = "Steven"

[void]$ws.Range("A27").Select()
